{"js": "// The <id>...</id> paragraph is currently split across three runs:\n//   \"<id>\"  (Courier New / #7F6000 / 9pt)\n//   \"p167r_1\"  (default run formatting)\n//   \"</id>\" (Courier New / #7F6000 / 9pt)\n// The edit merges them into a single run \"<id>p167r_1</id>\" that keeps the\n// first run's (Courier New) formatting.\nconst body = context.document.body;\n\n// Find the exact text span that covers all three runs.\nconst results = body.search(\"<id>p167r_1</id>\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text '<id>p167r_1</id>' not found in document body.\");\n}\n\nconst target = results.items[0];\n\n// Re-inserting the same text over the whole (multi-run) range collapses it\n// into a single run that inherits the formatting from the start of the\n// range (the first, Courier New run), exactly mirroring the OOXML diff.\ntarget.insertText(target.text, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The \"<id>p167r_1</id>\" paragraph is currently split across three runs:\n#   \"<id>\"     (Courier New / #7f6000 / 9pt)\n#   \"p167r_1\"  (default run formatting)\n#   \"</id>\"    (Courier New / #7f6000 / 9pt)\n# Use Find/Replace (not plain Range.Text assignment, which is a no-op when\n# the replacement text is identical to the existing text) to collapse the\n# three runs into a single run \"<id>p167r_1</id>\" that keeps the formatting\n# of the text at the start of the match (the first, Courier New run).\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"<id>p167r_1</id>\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"<id>p167r_1</id>\"\n\n# wdFindContinue = 1 (don't wrap past the end), wdReplaceOne = 1 (single match)\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 1) | Out-Null\n"}
